# Se procesan de nuevo los datos con las nuevas dimensiones curadas
#
# Row 2 (dimension names) / Row 3 (dimension kind) / Row 4 (dimension type)
# for column E ("provincia") are re-curated:
#   E2: sdmx-dimension:refArea      -> iaest-measure:provincia
#   E3: dim                         -> medida
#   E4: URI-Provincia                -> xsd:int

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Informe-02-020005-A-TC")

$ws.Range("E2").Value = "iaest-measure:provincia"
$ws.Range("E3").Value = "medida"
$ws.Range("E4").Value = "xsd:int"
